# "dev of account system"
# Append 17 new alliance/account-system related error rows (codes 623-639)
# to the "errors" sheet, right after the existing data (which ends at row 123).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: STR_key, INT_code, STR_message
$newErrors = @(
    @("allianceArchonCanNotQuitAlliance", 623, "别逗了,仅当联盟成员为空时,盟主才能退出联盟"),
    @("allianceInFightStatusCanNotQuitAlliance", 624, "联盟正在战争准备期或战争期,不能退出联盟"),
    @("allianceDoNotAllowJoinDirectly", 625, "联盟不允许直接加入"),
    @("joinAllianceRequestIsFull", 626, "联盟申请已满,请撤消部分申请后再来申请"),
    @("joinTheAllianceRequestAlreadySend", 627, "对此联盟的申请已发出,请耐心等候审核"),
    @("allianceJoinRequestMessagesIsFull", 628, "此联盟的申请信息已满,请等候其处理后再进行申请"),
    @("joinAllianceRequestNotExist", 629, "联盟申请事件不存在"),
    @("playerCancelTheJoinRequestToTheAlliance", 630, "玩家已经取消对此联盟的申请"),
    @("inviteRequestAlreadySend", 631, "此玩家已被邀请加入我方联盟,请等候其处理"),
    @("inviteRequestMessageIsFullForThisPlayer", 632, "此玩家的邀请信息已满,请等候其处理后再进行邀请"),
    @("allianceInviteEventNotExist", 633, "联盟邀请事件不存在"),
    @("playerAlreadyTheAllianceArchon", 634, "玩家已经是盟主了"),
    @("onlyAllianceArchonMoreThanSevenDaysNotOnLinePlayerCanBuyArchonTitle", 635, "盟主连续7天不登陆时才能购买盟主职位"),
    @("speedupRequestAlreadySendForThisEvent", 636, "此事件已经发送了加速请求"),
    @("allianceHelpEventNotExist", 637, "帮助事件不存在"),
    @("canNotHelpSelfSpeedup", 638, "不能帮助自己加速建造"),
    @("youAlreadyHelpedTheEvent", 639, "您已经帮助过此事件了")
)

$startRow = 124
$r = $startRow
foreach ($entry in $newErrors) {
    $ws.Range("A$r").Value = $entry[0]
    $ws.Range("B$r").Value = $entry[1]
    $ws.Range("C$r").Value = $entry[2]
    $r = $r + 1
}

$lastRow = $r - 1
$nextRow = $r

# Match the author's final selection/scroll state: cursor parked one row
# below the newly-added block, with the view scrolled down to it.
$ws.Range("A$nextRow").Select()
